$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 256.55554
$ws.Range("I4").Value = 263.625
$ws.Range("K4").Value = 263.625
$ws.Range("M4").Value = -149.625
$ws.Range("H43").Value = 4855.4443
$ws.Range("I43").Value = 4814.2856
$ws.Range("K43").Value = 4814.2856
$ws.Range("M43").Value = -4745.2856
$ws.Range("H96").Value = 45457900
$ws.Range("I96").Value = 71431660
$ws.Range("K96").Value = 214294980
$ws.Range("M96").Value = -214293607
$ws.Range("H100").Value = 489.66666
$ws.Range("I100").Value = 489.66666
$ws.Range("K100").Value = 489.66666
$ws.Range("M100").Value = 51.33334000000002
$ws.Range("H132").Value = 2672.1853
$ws.Range("I132").Value = 2006.25
$ws.Range("K132").Value = 6018.75
$ws.Range("M132").Value = -3488.75
$ws.Range("H135").Value = 2095.5334
$ws.Range("I135").Value = 1647.9
$ws.Range("J135").Value = 2990.8
$ws.Range("K135").Value = 14831.1
$ws.Range("L135").Value = 26917.2
$ws.Range("M135").Value = -12296.1
$ws.Range("N135").Value = -31987.2
$ws.Range("H137").Value = 6333.143
$ws.Range("I137").Value = 6220
$ws.Range("J137").Value = 6436
$ws.Range("K137").Value = 18660
$ws.Range("L137").Value = 19308
$ws.Range("M137").Value = -16110
$ws.Range("N137").Value = -24408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 7110.9375
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 7110.9375
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H97").Value = 709.25
$ws.Range("I97").Value = 695.6667
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 695.6667
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = -199.6667
$ws.Range("N97").Value = -1742
$ws.Range("H110").Value = 4465220.5
$ws.Range("I110").Value = 5495271.5
$ws.Range("J110").Value = 1666
$ws.Range("K110").Value = 5495271.5
$ws.Range("L110").Value = 1666
$ws.Range("M110").Value = -5493226.5
$ws.Range("N110").Value = -5756
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H94").Value = 1310.5333
$ws.Range("I94").Value = 1332
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 1332
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = -881
$ws.Range("N94").Value = -1912

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 6000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 6000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 6000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -6278
$ws.Range("H39").Value = 1412.375
$ws.Range("I39").Value = 1412.375
$ws.Range("K39").Value = 1412.375
$ws.Range("M39").Value = -1021.375
$ws.Range("H49").Value = 1412.375
$ws.Range("I49").Value = 1412.375
$ws.Range("K49").Value = 1412.375
$ws.Range("M49").Value = -1230.375
$ws.Range("H62").Value = 31067.312
$ws.Range("I62").Value = 5860.8335
$ws.Range("J62").Value = 106686.75
$ws.Range("K62").Value = 5860.8335
$ws.Range("L62").Value = 106686.75
$ws.Range("M62").Value = -5236.8335
$ws.Range("N62").Value = -107934.75
$ws.Range("H65").Value = 31067.312
$ws.Range("I65").Value = 5860.8335
$ws.Range("J65").Value = 106686.75
$ws.Range("K65").Value = 29304.1675
$ws.Range("L65").Value = 533433.75
$ws.Range("M65").Value = -26184.1675
$ws.Range("N65").Value = -539673.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 336.77777
$ws.Range("I13").Value = 128.625
$ws.Range("J13").Value = 2002
$ws.Range("K13").Value = 385.875
$ws.Range("L13").Value = 6006
$ws.Range("M13").Value = -217.875
$ws.Range("N13").Value = -6342
$ws.Range("H50").Value = 2411
$ws.Range("I50").Value = 1505
$ws.Range("J50").Value = 2637.5
$ws.Range("K50").Value = 4515
$ws.Range("L50").Value = 7912.5
$ws.Range("M50").Value = -4034
$ws.Range("N50").Value = -8874.5
$ws.Range("H53").Value = 2411
$ws.Range("I53").Value = 1505
$ws.Range("J53").Value = 2637.5
$ws.Range("K53").Value = 4515
$ws.Range("L53").Value = 7912.5
$ws.Range("M53").Value = -4034
$ws.Range("N53").Value = -8874.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9991.84
$ws.Range("I70").Value = 8380.549000000001
$ws.Range("J70").Value = 12620.789
$ws.Range("K70").Value = 8380.549000000001
$ws.Range("L70").Value = 12620.789
$ws.Range("M70").Value = -8110.549000000001
$ws.Range("N70").Value = -13160.789
$ws.Range("H73").Value = 9991.84
$ws.Range("I73").Value = 8380.549000000001
$ws.Range("J73").Value = 12620.789
$ws.Range("K73").Value = 8380.549000000001
$ws.Range("L73").Value = 12620.789
$ws.Range("M73").Value = -7444.549000000001
$ws.Range("N73").Value = -14492.789
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5428628
$ws.Range("I2").Value = 5000075
$ws.Range("J2").Value = 6000032.5
$ws.Range("K2").Value = 5000075
$ws.Range("L2").Value = 6000032.5
$ws.Range("M2").Value = -4999963
$ws.Range("N2").Value = -6000256.5
$ws.Range("H7").Value = 20986.625
$ws.Range("J7").Value = 20724.5
$ws.Range("L7").Value = 20724.5
$ws.Range("N7").Value = -20948.5
$ws.Range("H16").Value = 12077.941
$ws.Range("I16").Value = 13303.083
$ws.Range("J16").Value = 9137.6
$ws.Range("K16").Value = 13303.083
$ws.Range("L16").Value = 9137.6
$ws.Range("M16").Value = -13133.083
$ws.Range("N16").Value = -9477.6
$ws.Range("H46").Value = 6345.1665
$ws.Range("I46").Value = 2533.3333
$ws.Range("K46").Value = 2533.3333
$ws.Range("M46").Value = -2345.3333
$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("K48").Value = 1000
$ws.Range("M48").Value = -339
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180
$ws.Range("H116").Value = 10000
$ws.Range("J116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("N116").Value = -19178
$ws.Range("H126").Value = 20986.625
$ws.Range("J126").Value = 20724.5
$ws.Range("L126").Value = 62173.5
$ws.Range("N126").Value = -67113.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 399.4
$ws.Range("I13").Value = 200
$ws.Range("J13").Value = 449.25
$ws.Range("K13").Value = 200
$ws.Range("L13").Value = 449.25
$ws.Range("M13").Value = -60
$ws.Range("N13").Value = -729.25
$ws.Range("H17").Value = 5039
$ws.Range("I17").Value = 5048.75
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 5048.75
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = -4876.75
$ws.Range("N17").Value = -5344
$ws.Range("H122").Value = 2903.5715
$ws.Range("I122").Value = 3119.25
$ws.Range("J122").Value = 2616
$ws.Range("K122").Value = 9357.75
$ws.Range("L122").Value = 7848
$ws.Range("M122").Value = -6907.75
$ws.Range("N122").Value = -12748
